$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Ocupar mesa"
$ws.Range("B2:D2").Select()
